# This script rewrites the "output_data_dictionary" worksheet from a "wide"
# layout (separate City data / Grid data description columns) into a "long"
# layout (Category / Description / Variable / Scale columns), matching the
# updated data dictionary content and formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target content, row by row: Category | Description | Variable | Scale
$data = @(
    @('Category', 'Description', 'Variable', 'Scale'),
    @('Study region information', 'Continent', 'Continent', 'city'),
    @('Study region information', 'Country', 'Country', 'city'),
    @('Study region information', '2-letter country code', 'ISO 3166-1 alpha-2', 'city'),
    @('Study region information', 'Study region', 'City', 'city, grid'),
    @('Derived study region statistics', 'Area  (km²; accounting for urban restrictions, if applied)', 'Area (sqkm)', 'city, grid'),
    @('Derived study region statistics', 'Population estimate, as per configured population data source', 'Population estimate', 'city, grid'),
    @('Derived study region statistics', 'Population per km²', 'Population per sqkm', 'city, grid'),
    @('Derived study region statistics', 'Intersection count (following consolidation based on intersection tolerance parameter in region configuration)', 'Intersections', 'city, grid'),
    @('Derived study region statistics', 'Intersections per km²', 'Intersections per sqkm', 'city, grid'),
    @('Linked covariates', 'Total emission of CO 2 from the transport sector, using non-short-cycle-organic fuels in 2015', 'E_EC2E_T15', 'city'),
    @('Linked covariates', 'Total emission of CO 2 from the energy sector, using short-cycle-organic fuels in 2015', 'E_EC2O_T15', 'city'),
    @('Linked covariates', 'Total emission of PM 2.5 from the transport sector in 2015', 'E_EPM2_T15', 'city'),
    @('Linked covariates', 'Total concertation of PM 2.5 for reference epoch 2014', 'E_CPM2_T14', 'city'),
    @('Analytical statistic', 'Sample points used in this analysis (generated along pedestrian network for populated grid areas)', 'urban_sample_point_count', 'city, grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  fresh food market / supermarket (source: OpenStreetMap or custom)', 'access_500m_fresh_food_market_score', 'grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  convenience store (source: OpenStreetMap or custom)', 'access_500m_convenience_score', 'grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  public transport (source: OpenStreetMap or custom)', 'access_500m_pt_osm_any_score', 'grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  any public open space (source: OpenStreetMap)', 'access_500m_public_open_space_any_score', 'grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  public open space larger than 1.5 hectares (source: OpenStreetMap)', 'access_500m_public_open_space_large_score', 'grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  public transport (source: GTFS)', 'access_500m_pt_gtfs_any_score', 'grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  public transport with average daytime weekday service frequency of 30 minutes or better (source: GTFS)', 'access_500m_pt_gtfs_freq_30_score', 'grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  public transport with average daytime weekday service frequency of 20 minutes or better (source: GTFS)', 'access_500m_pt_gtfs_freq_20_score', 'grid'),
    @('Indicator estimates', 'Score (/1) for access within 500 m to a  any public transport stop (source: GTFS or OpenStreetMap/custom)', 'access_500m_pt_any_score', 'grid'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  fresh food market / supermarket (source: OpenStreetMap or custom)', 'pop_pct_access_500m_fresh_food_market_score', 'city'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  convenience store (source: OpenStreetMap or custom)', 'pop_pct_access_500m_convenience_score', 'city'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  public transport (source: OpenStreetMap or custom)', 'pop_pct_access_500m_pt_osm_any_score', 'city'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  any public open space (source: OpenStreetMap)', 'pop_pct_access_500m_public_open_space_any_score', 'city'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  public open space larger than 1.5 hectares (source: OpenStreetMap)', 'pop_pct_access_500m_public_open_space_large_score', 'city'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  public transport (source: GTFS)', 'pop_pct_access_500m_pt_gtfs_any_score', 'city'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  public transport with average daytime weekday service frequency of 30 minutes or better (source: GTFS)', 'pop_pct_access_500m_pt_gtfs_freq_30_score', 'city'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  public transport with average daytime weekday service frequency of 20 minutes or better (source: GTFS)', 'pop_pct_access_500m_pt_gtfs_freq_20_score', 'city'),
    @('Indicator estimates', 'Percentage of population with access within 500 m to a  any public transport stop (source: GTFS or OpenStreetMap/custom)', 'pop_pct_access_500m_pt_any_score', 'city'),
    @('Indicator estimates', 'Average walkable neighbourhood poulation density (population weighted) ', 'pop_nh_pop_density', 'city'),
    @('Indicator estimates', 'Average walkable neighbourhood intersection density (population weighted) ', 'pop_nh_intersection_density', 'city'),
    @('Indicator estimates', 'Average daily living score (population weighted)', 'pop_daily_living', 'city'),
    @('Indicator estimates', 'Average walkability (population weighted) ', 'pop_walkability', 'city, grid'),
    @('Indicator estimates', 'Average walkable neighbourhood poulation density ', 'local_nh_population_density', 'city, grid'),
    @('Indicator estimates', 'Average walkable neighbourhood intersection density ', 'local_nh_intersection_density', 'city, grid'),
    @('Indicator estimates', 'Average daily living score ', 'local_daily_living', 'city, grid'),
    @('Indicator estimates', 'Average walkability ', 'local_walkability', 'city, grid')
)

$rowCount = $data.Length
$colCount = 4

# Clear any pre-existing content/formatting in the sheet's used range so no
# stale values, borders or styles are left behind from the previous layout.
$ws.Cells.Clear()

# Write all cell values in one shot via a 2D array assigned to a Range.
$lastRow = $rowCount
$lastCol = $colCount
$targetRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item($lastRow, $lastCol))

$arr = New-Object 'object[,]' $rowCount, $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}
$targetRange.Value = $arr

# ---- Formatting ----
$xlContinuous = 1
$xlHairline = 1
$xlThin = 2
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeLeft = 7
$xlEdgeRight = 10
$xlVAlignTop = -4160
$xlHAlignLeft = -4131
$xlLineStyleNone = -4142

function Set-RowFormat($rowIndex, $isHeader, $topWeight, $bottomWeight) {
    $rowRange = $ws.Range($ws.Cells.Item($rowIndex,1), $ws.Cells.Item($rowIndex,4))
    $rowRange.Font.Name = "Calibri"
    $rowRange.Font.Size = 11
    $rowRange.Font.Color = 0
    $rowRange.Font.Bold = $isHeader
    $rowRange.VerticalAlignment = $xlVAlignTop
    $rowRange.WrapText = $false
    $rowRange.HorizontalAlignment = -4131 # general, overridden below for col B

    foreach ($c in 1..4) {
        $cell = $ws.Cells.Item($rowIndex, $c)
        $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone
        $cell.Borders.Item($xlEdgeRight).LineStyle = $xlLineStyleNone
        $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
        $cell.Borders.Item($xlEdgeTop).Weight = $topWeight
        $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
        $cell.Borders.Item($xlEdgeBottom).Weight = $bottomWeight
    }

    # Column B (Description) is left-aligned, top-aligned and wraps; columns
    # A/C/D only specify vertical alignment (top), general horizontal, no wrap.
    $colB = $ws.Cells.Item($rowIndex, 2)
    $colB.HorizontalAlignment = $xlHAlignLeft
    $colB.VerticalAlignment = $xlVAlignTop
    $colB.WrapText = $true

    foreach ($c in @(1,3,4)) {
        $cell = $ws.Cells.Item($rowIndex, $c)
        $cell.HorizontalAlignment = -4142 # xlGeneral - no explicit horizontal alignment stored
        $cell.VerticalAlignment = $xlVAlignTop
        $cell.WrapText = $false
    }
}

# Row 1: header - bold font, bottom border thin only (top none)
Set-RowFormat 1 $true $xlLineStyleNone $xlThin
foreach ($c in 1..4) {
    $ws.Cells.Item(1, $c).Borders.Item($xlEdgeTop).LineStyle = $xlLineStyleNone
}

# Row 2: first data row - top thin, bottom hair
Set-RowFormat 2 $false $xlThin $xlHairline

# Rows 3..(rowCount-1): middle data rows - top hair, bottom hair
for ($r = 3; $r -lt $rowCount; $r++) {
    Set-RowFormat $r $false $xlHairline $xlHairline
}

# Last row: top hair, bottom thin
Set-RowFormat $rowCount $false $xlHairline $xlThin

# ---- Sheet view state ----
$ws.Range("A39").Select()
